$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H40").Value = 6434.75
$ws.Range("I40").Value = 5829.2856
$ws.Range("J40").Value = 7282.4
$ws.Range("K40").Value = 5829.2856
$ws.Range("L40").Value = 7282.4
$ws.Range("M40").Value = -5654.2856
$ws.Range("N40").Value = -7632.4

$ws.Range("H44").Value = 49999
$ws.Range("I44").Value = 49999
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 49999
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -49537
$ws.Range("N44").ClearContents()

$ws.Range("H51").Value = 9333.333000000001
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 9333.333000000001
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9333.333000000001
$ws.Range("N51").Value = -10301.333

$ws.Range("H62").Value = 7946753
$ws.Range("I62").Value = 11001514
$ws.Range("J62").Value = 4375.8
$ws.Range("K62").Value = 11001514
$ws.Range("L62").Value = 4375.8
$ws.Range("M62").Value = -11000890
$ws.Range("N62").Value = -5623.8

$ws.Range("H65").Value = 7946753
$ws.Range("I65").Value = 11001514
$ws.Range("J65").Value = 4375.8
$ws.Range("K65").Value = 55007570
$ws.Range("L65").Value = 21879
$ws.Range("M65").Value = -55004450
$ws.Range("N65").Value = -28119

$ws.Range("H135").Value = 7200.8335
$ws.Range("I135").Value = 4773
$ws.Range("J135").Value = 10599.8
$ws.Range("K135").Value = 42957
$ws.Range("L135").Value = 95398.2
$ws.Range("M135").Value = -40422
$ws.Range("N135").Value = -100468.2

$ws.Range("H138").Value = 4114.6943
$ws.Range("I138").Value = 1284
$ws.Range("J138").Value = 5916.0454
$ws.Range("K138").Value = 3852
$ws.Range("L138").Value = 17748.1362
$ws.Range("M138").Value = 1288
$ws.Range("N138").Value = -28028.1362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 54120
$ws.Range("I43").Value = 42600
$ws.Range("J43").Value = 57000
$ws.Range("K43").Value = 42600
$ws.Range("L43").Value = 57000
$ws.Range("M43").Value = -42287
$ws.Range("N43").Value = -57626

$ws.Range("H63").Value = 3621.2778
$ws.Range("I63").Value = 1677.1428
$ws.Range("J63").Value = 4858.4546
$ws.Range("K63").Value = 1677.1428
$ws.Range("L63").Value = 4858.4546
$ws.Range("M63").Value = -991.1428000000001
$ws.Range("N63").Value = -6230.4546

$ws.Range("H66").Value = 3621.2778
$ws.Range("I66").Value = 1677.1428
$ws.Range("J66").Value = 4858.4546
$ws.Range("K66").Value = 8385.714
$ws.Range("L66").Value = 24292.273
$ws.Range("M66").Value = -4953.714
$ws.Range("N66").Value = -31156.273

$ws.Range("H122").Value = 2605.5667
$ws.Range("I122").Value = 2069
$ws.Range("J122").Value = 3857.5557
$ws.Range("K122").Value = 6207
$ws.Range("L122").Value = 11572.6671
$ws.Range("M122").Value = -3757
$ws.Range("N122").Value = -16472.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11800.48
$ws.Range("I86").Value = 14747
$ws.Range("J86").Value = 9836.134
$ws.Range("K86").Value = 14747
$ws.Range("L86").Value = 9836.134
$ws.Range("M86").Value = -13624
$ws.Range("N86").Value = -12082.134

$ws.Range("H89").Value = 11800.48
$ws.Range("I89").Value = 14747
$ws.Range("J89").Value = 9836.134
$ws.Range("K89").Value = 73735
$ws.Range("L89").Value = 49180.67
$ws.Range("M89").Value = -68119
$ws.Range("N89").Value = -60412.67

$ws.Range("H107").Value = 8335537.5
$ws.Range("I107").Value = 10002256
$ws.Range("J107").Value = 1947
$ws.Range("K107").Value = 10002256
$ws.Range("L107").Value = 1947
$ws.Range("M107").Value = -10000336
$ws.Range("N107").Value = -5787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2658.762
$ws.Range("I31").Value = 2015.2858
$ws.Range("J31").Value = 2980.5
$ws.Range("K31").Value = 2015.2858
$ws.Range("L31").Value = 2980.5
$ws.Range("M31").Value = -1720.2858
$ws.Range("N31").Value = -3570.5

$ws.Range("H34").Value = 2658.762
$ws.Range("I34").Value = 2015.2858
$ws.Range("J34").Value = 2980.5
$ws.Range("K34").Value = 2015.2858
$ws.Range("L34").Value = 2980.5
$ws.Range("M34").Value = -1813.2858
$ws.Range("N34").Value = -3384.5

$ws.Range("H35").Value = 1739.75
$ws.Range("I35").Value = 1739.75
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1739.75
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1445.75

$ws.Range("H122").Value = 3682.2083
$ws.Range("I122").Value = 3504.8125
$ws.Range("J122").Value = 4037
$ws.Range("K122").Value = 10514.4375
$ws.Range("L122").Value = 12111
$ws.Range("M122").Value = -8064.4375
$ws.Range("N122").Value = -17011

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4450379.5
$ws.Range("I4").Value = 4584372
$ws.Range("J4").Value = 1234567
$ws.Range("K4").Value = 13753116
$ws.Range("L4").Value = 3703701
$ws.Range("M4").Value = -13753004
$ws.Range("N4").Value = -3703925

$ws.Range("H56").Value = 8774.706
$ws.Range("I56").Value = 8774.706
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 8774.706
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8244.706

$ws.Range("H101").Value = 27500
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 27500
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 82500
$ws.Range("N101").Value = -87368

$ws.Range("H131").Value = 32100826
$ws.Range("I131").Value = 76190920
$ws.Range("J131").Value = 16669294
$ws.Range("K131").Value = 228572760
$ws.Range("L131").Value = 50007882
$ws.Range("M131").Value = -228567720
$ws.Range("N131").Value = -50017962

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5488
$ws.Range("I16").Value = 1573.5
$ws.Range("J16").Value = 9961.714
$ws.Range("K16").Value = 1573.5
$ws.Range("L16").Value = 9961.714
$ws.Range("M16").Value = -1403.5
$ws.Range("N16").Value = -10301.714

$ws.Range("H82").Value = 3039.7693
$ws.Range("I82").Value = 1949.8334
$ws.Range("J82").Value = 3974
$ws.Range("K82").Value = 1949.8334
$ws.Range("L82").Value = 3974
$ws.Range("M82").Value = -1588.8334
$ws.Range("N82").Value = -4696

$ws.Range("H85").Value = 3039.7693
$ws.Range("I85").Value = 1949.8334
$ws.Range("J85").Value = 3974
$ws.Range("K85").Value = 1949.8334
$ws.Range("L85").Value = 3974
$ws.Range("M85").Value = -701.8334
$ws.Range("N85").Value = -6470

$ws.Range("H101").Value = 34802
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 34802
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 34802
$ws.Range("N101").Value = -41292

$ws.Range("H132").Value = 2341.2917
$ws.Range("I132").Value = 2033.6818
$ws.Range("J132").Value = 5725
$ws.Range("K132").Value = 6101.0454
$ws.Range("L132").Value = 17175
$ws.Range("M132").Value = -3571.0454
$ws.Range("N132").Value = -22235

$ws.Range("H136").Value = 8200010
$ws.Range("I136").Value = 21739964
$ws.Range("J136").Value = 4774
$ws.Range("K136").Value = 65219892
$ws.Range("L136").Value = 14322
$ws.Range("M136").Value = -65217342
$ws.Range("N136").Value = -19422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1715.75
$ws.Range("I81").Value = 1581.8889
$ws.Range("J81").Value = 2117.3333
$ws.Range("K81").Value = 3163.7778
$ws.Range("L81").Value = 4234.6666
$ws.Range("M81").Value = -2102.7778
$ws.Range("N81").Value = -6356.6666

$ws.Range("H84").Value = 1715.75
$ws.Range("I84").Value = 1581.8889
$ws.Range("J84").Value = 2117.3333
$ws.Range("K84").Value = 15818.889
$ws.Range("L84").Value = 21173.333
$ws.Range("M84").Value = -10514.889
$ws.Range("N84").Value = -31781.333

$ws.Range("H95").Value = 27983.25
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 27983.25
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 27983.25
$ws.Range("N95").Value = -33475.25

$ws.Range("H100").Value = 1347.1765
$ws.Range("I100").Value = 939.1429000000001
$ws.Range("J100").Value = 1632.8
$ws.Range("K100").Value = 1878.2858
$ws.Range("L100").Value = 3265.6
$ws.Range("M100").Value = -1337.2858
$ws.Range("N100").Value = -4347.6

$ws.Range("H122").Value = 2574.5386
$ws.Range("I122").Value = 1997.1818
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 5991.5454
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -3541.5454
$ws.Range("N122").Value = -22150

$ws.Range("H136").Value = 125259870
$ws.Range("I136").Value = 167003250
$ws.Range("J136").Value = 29750
$ws.Range("K136").Value = 501009750
$ws.Range("L136").Value = 89250
$ws.Range("M136").Value = -501007200
$ws.Range("N136").Value = -94350
